$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The comment in E7 ("saturating_two_bits" row) duplicated/contained a mistaken
# variant of the comment already used in E6 ("saturating_one_bit" row). Fix the
# mistake by reusing the correct text (this also removes the now-unused unique
# shared string for the erroneous comment).
$ws.Range("E7").Value = $ws.Range("E6").Value()

# Update the active cell selection to C8, matching the saved view state.
[void]$ws.Range("C8").Select()
